$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = -21.58930000000001
$ws.Range("A27").Value = -21.8019
$ws.Range("A32").Value = -21.2682
$ws.Range("A36").Value = -20.2051
$ws.Range("A38").Value = -19.6736
$ws.Range("A46").Value = -21.86680000000001
$ws.Range("A54").Value = -21.70200000000001
$ws.Range("A55").Value = -22.56920000000002
$ws.Range("A56").Value = -22.2008
$ws.Range("A67").Value = -21.49759999999997
$ws.Range("A69").Value = -21.56109999999997
$ws.Range("A72").Value = -21.62949999999999
$ws.Range("A83").Value = -21.71359999999999
$ws.Range("A86").Value = -22.13210000000002
$ws.Range("A91").Value = -21.45170000000001
$ws.Range("A93").Value = -21.10879999999999
$ws.Range("A99").Value = -20.15439999999999
